$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple two-row swaps of columns B:G (stock report rows were re-sequenced)
$row142 = $ws.Range("B142:G142").Value()
$row143 = $ws.Range("B143:G143").Value()
$ws.Range("B142:G142").Value = $row143
$ws.Range("B143:G143").Value = $row142

$row305 = $ws.Range("B305:G305").Value()
$row306 = $ws.Range("B306:G306").Value()
$ws.Range("B305:G305").Value = $row306
$ws.Range("B306:G306").Value = $row305

$row338 = $ws.Range("B338:G338").Value()
$row339 = $ws.Range("B339:G339").Value()
$ws.Range("B338:G338").Value = $row339
$ws.Range("B339:G339").Value = $row338

$row342 = $ws.Range("B342:G342").Value()
$row344 = $ws.Range("B344:G344").Value()
$ws.Range("B342:G342").Value = $row344
$ws.Range("B344:G344").Value = $row342

$row347 = $ws.Range("B347:G347").Value()
$row348 = $ws.Range("B348:G348").Value()
$ws.Range("B347:G347").Value = $row348
$ws.Range("B348:G348").Value = $row347

$row364 = $ws.Range("B364:G364").Value()
$row365 = $ws.Range("B365:G365").Value()
$ws.Range("B364:G364").Value = $row365
$ws.Range("B365:G365").Value = $row364

$row367 = $ws.Range("B367:G367").Value()
$row368 = $ws.Range("B368:G368").Value()
$ws.Range("B367:G367").Value = $row368
$ws.Range("B368:G368").Value = $row367

$row374 = $ws.Range("B374:G374").Value()
$row375 = $ws.Range("B375:G375").Value()
$ws.Range("B374:G374").Value = $row375
$ws.Range("B375:G375").Value = $row374

$row381 = $ws.Range("B381:G381").Value()
$row382 = $ws.Range("B382:G382").Value()
$ws.Range("B381:G381").Value = $row382
$ws.Range("B382:G382").Value = $row381

$row392 = $ws.Range("B392:G392").Value()
$row393 = $ws.Range("B393:G393").Value()
$ws.Range("B392:G392").Value = $row393
$ws.Range("B393:G393").Value = $row392

$row413 = $ws.Range("B413:G413").Value()
$row414 = $ws.Range("B414:G414").Value()
$ws.Range("B413:G413").Value = $row414
$ws.Range("B414:G414").Value = $row413

$row423 = $ws.Range("B423:G423").Value()
$row424 = $ws.Range("B424:G424").Value()
$ws.Range("B423:G423").Value = $row424
$ws.Range("B424:G424").Value = $row423

$row528 = $ws.Range("B528:G528").Value()
$row529 = $ws.Range("B529:G529").Value()
$ws.Range("B528:G528").Value = $row529
$ws.Range("B529:G529").Value = $row528

$row578 = $ws.Range("B578:G578").Value()
$row579 = $ws.Range("B579:G579").Value()
$ws.Range("B578:G578").Value = $row579
$ws.Range("B579:G579").Value = $row578

$row582 = $ws.Range("B582:G582").Value()
$row583 = $ws.Range("B583:G583").Value()
$ws.Range("B582:G582").Value = $row583
$ws.Range("B583:G583").Value = $row582

$row585 = $ws.Range("B585:G585").Value()
$row586 = $ws.Range("B586:G586").Value()
$ws.Range("B585:G585").Value = $row586
$ws.Range("B586:G586").Value = $row585

$row701 = $ws.Range("B701:G701").Value()
$row702 = $ws.Range("B702:G702").Value()
$ws.Range("B701:G701").Value = $row702
$ws.Range("B702:G702").Value = $row701

$row712 = $ws.Range("B712:G712").Value()
$row713 = $ws.Range("B713:G713").Value()
$ws.Range("B712:G712").Value = $row713
$ws.Range("B713:G713").Value = $row712

$row864 = $ws.Range("B864:G864").Value()
$row865 = $ws.Range("B865:G865").Value()
$ws.Range("B864:G864").Value = $row865
$ws.Range("B865:G865").Value = $row864

# 3-row rotation: row154 <- row155 <- row156 <- row154 (wrap)
$row154 = $ws.Range("B154:G154").Value()
$row155 = $ws.Range("B155:G155").Value()
$row156 = $ws.Range("B156:G156").Value()
$ws.Range("B154:G154").Value = $row155
$ws.Range("B155:G155").Value = $row156
$ws.Range("B156:G156").Value = $row154

# 3-row rotation: row308 <- row310, row309 <- row308, row310 <- row309 (wrap the other way)
$row308 = $ws.Range("B308:G308").Value()
$row309 = $ws.Range("B309:G309").Value()
$row310 = $ws.Range("B310:G310").Value()
$ws.Range("B308:G308").Value = $row310
$ws.Range("B309:G309").Value = $row308
$ws.Range("B310:G310").Value = $row309
